$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.125.81'
$ws.Cells.Item(2, 5).Value = '  -4.49%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.651.69'
$ws.Cells.Item(3, 5).Value = '  -3.71%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.013'
$ws.Cells.Item(4, 5).Value = '  +0.53%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.68'
$ws.Cells.Item(5, 5).Value = '  -4.08%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5113'
$ws.Cells.Item(6, 5).Value = '  -3.18%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.46%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2591'
$ws.Cells.Item(8, 5).Value = '  -2.16%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06431'
$ws.Cells.Item(9, 5).Value = '  -3.53%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.68'
$ws.Cells.Item(10, 5).Value = '  -5.26%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07785'
$ws.Cells.Item(11, 5).Value = '  +0.38%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.654.48'
$ws.Cells.Item(12, 5).Value = '  -3.41%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.274'
$ws.Cells.Item(13, 5).Value = '  -4.42%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.882.76'
$ws.Cells.Item(14, 5).Value = '  -3.49%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.5468'
$ws.Cells.Item(15, 5).Value = '  -5.69%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₅7986'
$ws.Cells.Item(16, 5).Value = '  -2.32%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.85'
$ws.Cells.Item(17, 5).Value = '  -5.83%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '26.167.88'
$ws.Cells.Item(18, 5).Value = '  -4.36%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.01%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '207.66'
$ws.Cells.Item(20, 5).Value = '  -5.44%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.391'
$ws.Cells.Item(21, 5).Value = '  -5.59%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.06'
$ws.Cells.Item(22, 5).Value = '  -3.44%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '6.035'
$ws.Cells.Item(23, 5).Value = '  -0.07%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.43%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.854'
$ws.Cells.Item(25, 5).Value = '  +7.45%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '144.19'
$ws.Cells.Item(26, 5).Value = '  -0.57%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.1168'
$ws.Cells.Item(27, 5).Value = '  -3.19%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.938'
$ws.Cells.Item(28, 5).Value = '  -3.98%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '15.77'
$ws.Cells.Item(29, 5).Value = '  -2.62%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -4.93%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.243'
$ws.Cells.Item(31, 5).Value = '  -3.94%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.336'
$ws.Cells.Item(32, 5).Value = '  -4.15%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.234'
$ws.Cells.Item(33, 5).Value = '  -4.76%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.543'
$ws.Cells.Item(34, 5).Value = '  -5.85%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.359'
$ws.Cells.Item(35, 5).Value = '  -1.70%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.711'
$ws.Cells.Item(36, 5).Value = '  -4.70%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.9162'
$ws.Cells.Item(37, 5).Value = '  -3.83%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.171.25'
$ws.Cells.Item(38, 5).Value = '  -0.94%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.5690'
$ws.Cells.Item(39, 5).Value = '  -3.10%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01580'
$ws.Cells.Item(40, 5).Value = '  -4.22%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'mCoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.579'
$ws.Cells.Item(41, 5).Value = '  +0.15%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'PaxDollar'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.011'
$ws.Cells.Item(42, 5).Value = '  +0.32%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.663'
$ws.Cells.Item(43, 5).Value = '  -2.80%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.8255'
$ws.Cells.Item(44, 5).Value = '  -1.63%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '100.24'
$ws.Cells.Item(45, 5).Value = '  -0.84%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '1.794.12'
$ws.Cells.Item(46, 5).Value = '  -3.43%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.0₈114'
$ws.Cells.Item(47, 5).Value = '  -2.92%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4562'
$ws.Cells.Item(48, 5).Value = '  +0.32%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.06%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '55.20'
$ws.Cells.Item(50, 5).Value = '  -4.05%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.842'
$ws.Cells.Item(51, 5).Value = '  -3.88%  '

